$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.731.42"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "2.022.21"
$ws.Range("E3").Value = "  -1.07%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'226.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.68%  "
$ws.Range("D6").Value = "'0.610"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.63%  "
$ws.Range("D7").Value = "'59.34"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.06%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.382"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("D10").Value = "'0.0804"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.04%  "
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "2.322.57"
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'14.47"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("D14").Value = "'20.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.61%  "
$ws.Range("D15").Value = "'0.747"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.41%  "
$ws.Range("D16").Value = "'5.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("D17").Value = "2.031.36"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").Value = "37.704.59"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").Value = "'6.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.86%  "
$ws.Range("D20").Value = "'69.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").Value = "0.0₃0819"
$ws.Range("E21").Value = "  -0.95%  "
$ws.Range("D22").Value = "'223.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "'2.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.44%  "
$ws.Range("D25").Value = "'2.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.59%  "
$ws.Range("D26").Value = "'165.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("D27").Value = "'9.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("D28").Value = "'0.129"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.11%  "
$ws.Range("D29").Value = "'18.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.56%  "
$ws.Range("D30").Value = "'1.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.10%  "
$ws.Range("D31").Value = "'0.119"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.06%  "
$ws.Range("D32").Value = "'4.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.86%  "
$ws.Range("E33").Value = "  +0.77%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.0599"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.90%  "
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "'4.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.16%  "
$ws.Range("E36").Value = "  +6.09%  "
$ws.Range("D37").Value = "'2.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.04%  "
$ws.Range("D38").Value = "'3.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.07%  "
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").Value = "1.531.20"
$ws.Range("E40").Value = "  +3.42%  "
$ws.Range("D41").Value = "'0.0215"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("D42").Value = "'95.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.61%  "
$ws.Range("D43").Value = "'16.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").Value = "'2.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.06%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "'0.0915"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.69%  "
$ws.Range("E46").Value = "  -1.47%  "
$ws.Range("B47").Value = "MXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D47").Value = "'2.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.63%  "
$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D48").Value = "'3.87"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.58%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "'0.998"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.72%  "
$ws.Range("D50").Value = "'7.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.83%  "
$ws.Range("D51").Value = "2.213.99"
$ws.Range("E51").Value = "  -0.79%  "
